# Re-pull / push updated data for the "dSF" (F) column on specific rows.
# Matches commit message: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    3  = 9
    5  = -8
    12 = -4
    13 = 3
    20 = 0
    21 = -6
    23 = -4
    26 = -3
    37 = -6
    39 = 6
    40 = -12
    41 = -5
    42 = 3
    43 = -6
    44 = -2
    45 = -4
    50 = -3
    51 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
